$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 31 (pushing "bloom" and everything below down by one),
# matching the formatting of the surrounding rows (default style, row height 13.4).
$ws.Rows.Item(31).Insert()

$ws.Range("A31").Value = "god rays"
$ws.Range("B31").Value = "God rays"

$ws.Range("A31:B31").RowHeight = 13.4

# The insert shifts row heights down; re-assert the explicit (non-default)
# heights for the rows that fell out of sync with the default.
$ws.Range("A103:B108").RowHeight = 12.8
$ws.Range("A155:B155").RowHeight = 12.8
